$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FNEC")

# Row 58 = "Short/Current Long Term Debt" (was all zeros, now holds the values
# that used to live in row 59)
$ws.Range("D58").Value = 700
$ws.Range("E58").Value = 700
$ws.Range("F58").Value = 700
$ws.Range("G58").Value = 600
$ws.Range("H58").Value = 600
$ws.Range("I58").Value = 600
$ws.Range("J58").Value = 500

# Row 59 = "Other Current Liabilities" (now "NA" for D:I, and 0 for J)
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "NA"
$ws.Range("H59").Value = "NA"
$ws.Range("I59").Value = "NA"
$ws.Range("J59").Value = 0

$wb.Save()
